# CoronaVirus US Lockdown Forecast - Apr 1 data update
# I35 becomes a hard-coded actual value (215003) instead of the forecast
# formula; everything downstream (I36:I49 and the J/K/L/M/N columns for
# rows 35-49) recalculates automatically off of that new anchor.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy I34's formatting (the "actual data" fill/number format) onto I35 so
# the cell matches the style used by the other hard-coded actuals (I32:I34)
# instead of the forecast-formula fill it currently has.
$ws.Range("I34").Copy()
$ws.Range("I35").PasteSpecial(-4122)

# Replace the forecast formula in I35 with the real reported value.
$ws.Range("I35").Value = 215003

# Move the active selection to I36, matching where the author left off.
$ws.Range("I36").Select()
